# Auto-generated Excel COM-interop script to update cryptos list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text while we write values, so that
# numeric-looking strings (e.g. "506.40", "1.00") are not silently
# converted into actual numbers and lose their original text formatting.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "54.745.79"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "2.280.83"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "506.40"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").Value = "128.93"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "2.302.94"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "0.0969"
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "0.346"
$ws.Range("E12").Value = "  +4.97%  "
$ws.Range("D13").Value = "4.93"
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("D14").Value = "23.44"
$ws.Range("E14").Value = "  +5.55%  "
$ws.Range("D15").Value = "2.689.73"
$ws.Range("D16").Value = "54.800.65"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "2.302.55"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "307.14"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").Value = "6.42"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "60.12"
$ws.Range("E24").Value = "  -2.60%  "
$ws.Range("D25").Value = "0.994"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").Value = "7.45"
$ws.Range("E27").Value = "  +3.36%  "
$ws.Range("D28").Value = "170.96"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "0.0₃0705"
$ws.Range("E29").Value = "  +3.33%  "
$ws.Range("D30").Value = "6.07"
$ws.Range("E30").Value = "  +3.52%  "
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").Value = "0.995"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "0.922"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").Value = "3.78"
$ws.Range("E38").Value = "  +2.37%  "
$ws.Range("D39").Value = "36.38"
$ws.Range("D40").Value = "0.376"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "5.06"
$ws.Range("E42").Value = "  +5.99%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.41"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "126.07"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "0.0498"
$ws.Range("E45").Value = "  +2.32%  "
$ws.Range("D46").Value = "248.99"
$ws.Range("E46").Value = "  +4.45%  "
$ws.Range("D47").Value = "0.0905"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").Value = "0.551"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("E51").Value = "  +0.51%  "

# Restore column D to the default "Normal" style so no stray number
# formatting is left behind on cells that did not have one originally.
$dRange.Style = "Normal"

